$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Enterprises density (per 1000 people)" indicator row (originally the
# third row of the block, A12:D12) moves up to become the first row of the
# block (A10:D10). "Employment (% of total)" and "Enterprises (absolute #)"
# each shift down by one row. "Employment (absolute #)" (row 13) and
# "Enterprises (% of total)" (row 14) are unaffected.
#
# Cell-by-cell Copy/PasteSpecial (instead of re-typing the .Value) keeps the
# original cell type (text, since "33.6" etc. are stored as strings, not
# numbers) and the original cell style intact.

$xlPasteAll = -4104

# Stash the density row (A12:D12) in a scratch cell out of the printed area.
$ws.Range("A12:D12").Copy() | Out-Null
$ws.Range("A100:D100").PasteSpecial($xlPasteAll) | Out-Null

# Shift "Employment (% of total)" (row 10) down into row 11's old spot... 
# done from the bottom up so we never overwrite a row before it's been saved.
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial($xlPasteAll) | Out-Null

$ws.Range("A10:D10").Copy() | Out-Null
$ws.Range("A11:D11").PasteSpecial($xlPasteAll) | Out-Null

# Drop the density row into the now-vacated top slot.
$ws.Range("A100:D100").Copy() | Out-Null
$ws.Range("A10:D10").PasteSpecial($xlPasteAll) | Out-Null

# Clean up the scratch cells.
$ws.Range("A100:D100").Clear() | Out-Null

$excel.CutCopyMode = $false
